$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the comment text for Approach 2 qualitative evaluation (row 11, col B):
# "many different datasets" -> "for different datasets"
$ws.Range("B11").Value = "Used for qualitative evaluation of Approach 2 on both original and retrained models for different datasets. "

# Update the active selection to match the saved cursor position recorded in the file
$ws.Range("B12").Select()
